$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary block (rows 10-12) ---

# Row 10: No.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("B10").Value = 22
$ws.Range("D10").Value = 6
$ws.Range("E10").Value = 28

# Row 11: Marking
$ws.Range("A9").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Row 12: Total
$ws.Range("A9").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("B12").Value = 88
$ws.Range("E12").Value = "88/112"

# --- Remove the third answer block (columns G:H) entirely ---
$ws.Range("G15:H40").Clear()

# --- Remove the second answer block (columns D:E) except rows 16-18 ---
$ws.Range("D19:E40").Clear()

# Fill in D16:D18 (Student Ans for block 2) to match the Correct Ans in E16:E18
$ws.Range("B10").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = "Option A"

$ws.Range("B10").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D17").Value = "Option C"

$ws.Range("B10").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = "Option D"

# --- Fill the Student Ans column (A) in block 1 for the rows that now have an answer ---
$answers = @{
    16 = "Option A"
    17 = "Option D"
    18 = "Option B"
    19 = "Option C"
    20 = "Option B"
    21 = "Option C"
    22 = "Option D"
    28 = "Option D"
    29 = "Option D"
    30 = "Option B"
    31 = "Option D"
    33 = "Option D"
    34 = "Option B"
    35 = "Option D"
    36 = "Option A"
    37 = "Option A"
    38 = "Option A"
    39 = "Option D"
    40 = "Option D"
}

foreach ($r in $answers.Keys) {
    $ws.Range("B10").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)
    $ws.Range("A$r").Value = $answers[$r]
}
